# Update Name of Algo
# Updates the imputed values in column C for a set of rows in Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "C9"  = -10.36670000000001
    "C13" = -12.4265
    "C16" = -12.82250000000001
    "C18" = -11.23459999999999
    "C20" = -12.0743
    "C26" = -12.325
    "C27" = -12.45569999999999
    "C29" = -11.6773
    "C35" = -12.1235
    "C36" = -11.89320000000001
    "C45" = -13.97289999999998
    "C55" = -13.6319
    "C57" = -13.83429999999999
    "C69" = -10.9467
    "C76" = -12.6091
    "C78" = -13.0782
    "C82" = -11.8944
    "C83" = -14.09349999999999
    "C93" = -10.31449999999999
    "C97" = -12.3759
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
